$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.789.53"
$ws.Range("E2").Value = "  +2.87%  "

# Row 3
$ws.Range("D3").Value = "2.604.18"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.84%  "

# Row 9
$ws.Range("D9").Value = "2.628.23"
$ws.Range("E9").Value = "  +1.82%  "

# Row 10
$ws.Range("E10").Value = "  -2.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.74%  "

# Row 12
$ws.Range("E12").Value = "  -3.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.368"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.05%  "

# Row 14
$ws.Range("D14").Value = "3.069.28"
$ws.Range("E14").Value = "  +1.21%  "

# Row 15
$ws.Range("D15").Value = "60.776.60"
$ws.Range("E15").Value = "  +2.76%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.32%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.13%  "

# Row 18
$ws.Range("D18").Value = "2.617.09"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.61%  "

# Row 23
$ws.Range("E23").Value = "  +0.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.526"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.40%  "

# Row 25
$ws.Range("E25").Value = "  -0.37%  "

# Row 26
$ws.Range("E26").Value = "  -1.04%  "

# Row 27
$ws.Range("E27").Value = "  -0.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "0.0₃0794"
$ws.Range("E29").Value = "  +1.85%  "

# Row 30
$ws.Range("E30").Value = "  +8.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.961"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.74%  "

# Row 37
$ws.Range("E37").Value = "  +3.53%  "

# Row 38
$ws.Range("E38").Value = "  +7.32%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.850"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "297.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.27%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0241"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.67%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
